$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.299.23'
$ws.Range('E2').Value = '  -0.01%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.919.90'
$ws.Range('E3').Value = '  -0.34%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.86%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '485.20'
$ws.Range('E5').Value = '  +1.95%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.45'
$ws.Range('E6').Value = '  -0.19%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  +1.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('E8').Value = '  -0.23%  '

# Row 9
$ws.Range('E9').Value = '  +2.15%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  +2.65%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000348'
$ws.Range('E11').Value = '  +0.83%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.22'
$ws.Range('E12').Value = '  +1.63%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.76'
$ws.Range('E13').Value = '  +5.28%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.541.83'
$ws.Range('E14').Value = '  -1.68%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.900.35'
$ws.Range('E15').Value = '  -2.30%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.27'
$ws.Range('E16').Value = '  -2.09%  '

# Row 17
$ws.Range('E17').Value = '  -0.52%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.24'
$ws.Range('E18').Value = '  +2.77%  '

# Row 19
$ws.Range('E19').Value = '  +1.89%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.324.78'
$ws.Range('E20').Value = '  -0.84%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '431.96'
$ws.Range('E21').Value = '  -0.07%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.54'
$ws.Range('E22').Value = '  +8.09%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.12'
$ws.Range('E23').Value = '  +6.12%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.66'
$ws.Range('E24').Value = '  +3.07%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.69'
$ws.Range('E25').Value = '  +20.43%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.71'
$ws.Range('E26').Value = '  +2.26%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.20'
$ws.Range('E27').Value = '  +13.22%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.82'
$ws.Range('E28').Value = '  -0.52%  '

# Row 29
$ws.Range('E29').Value = '  -1.40%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '719.23'
$ws.Range('E30').Value = '  +0.39%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.77'
$ws.Range('E31').Value = '  +4.95%  '

# Row 32
$ws.Range('E32').Value = '  +3.44%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.93'
$ws.Range('E33').Value = '  +4.88%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.19'
$ws.Range('E34').Value = '  +17.20%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0891'
$ws.Range('E35').Value = '  +8.75%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '41.74'
$ws.Range('E36').Value = '  +0.45%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '61.12'
$ws.Range('E37').Value = '  -4.06%  '

# Row 38
$ws.Range('B38').Value = 'TheGraph'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.405'
$ws.Range('E38').Value = '  +23.31%  '

# Row 39
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.07'
$ws.Range('E39').Value = '  +19.07%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  +0.08%  '

# Row 41
$ws.Range('E41').Value = '  -2.93%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0493'
$ws.Range('E42').Value = '  +5.63%  '

# Row 43
$ws.Range('E43').Value = '  +4.46%  '

# Row 44
$ws.Range('E44').Value = '  +3.32%  '

# Row 45
$ws.Range('E45').Value = '  +1.80%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.36'
$ws.Range('E46').Value = '  +5.20%  '

# Row 47
$ws.Range('E47').Value = '  -0.58%  '

# Row 48
$ws.Range('E48').Value = '  +2.59%  '

# Row 49
$ws.Range('E49').Value = '  +0.70%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.32'
$ws.Range('E50').Value = '  -0.70%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0336'
$ws.Range('E51').Value = '  +32.38%  '
